$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 330.375
$ws.Range("I12").Value = 330.375
$ws.Range("K12").Value = 330.375
$ws.Range("M12").Value = -160.375
$ws.Range("H62").Value = 3176.889
$ws.Range("I62").Value = 3010.875
$ws.Range("J62").Value = 3309.7
$ws.Range("K62").Value = 3010.875
$ws.Range("L62").Value = 3309.7
$ws.Range("M62").Value = -2386.875
$ws.Range("N62").Value = -4557.7
$ws.Range("H63").Value = 41666.668
$ws.Range("J63").Value = 41666.668
$ws.Range("L63").Value = 41666.668
$ws.Range("N63").Value = -42914.668
$ws.Range("H65").Value = 3176.889
$ws.Range("I65").Value = 3010.875
$ws.Range("J65").Value = 3309.7
$ws.Range("K65").Value = 15054.375
$ws.Range("L65").Value = 16548.5
$ws.Range("M65").Value = -11934.375
$ws.Range("N65").Value = -22788.5
$ws.Range("H66").Value = 41666.668
$ws.Range("J66").Value = 41666.668
$ws.Range("L66").Value = 125000.004
$ws.Range("N66").Value = -131240.004
$ws.Range("H80").Value = 933.3333
$ws.Range("I80").Value = 820
$ws.Range("K80").Value = 2460
$ws.Range("M80").Value = -1462
$ws.Range("H83").Value = 933.3333
$ws.Range("I83").Value = 820
$ws.Range("K83").Value = 7380
$ws.Range("M83").Value = -2388
$ws.Range("H116").Value = 42634696
$ws.Range("I116").Value = 28973908
$ws.Range("K116").Value = 28973908
$ws.Range("M116").Value = -28970466
$ws.Range("H132").Value = 12374.52
$ws.Range("I132").Value = 5376.909
$ws.Range("J132").Value = 14251.927
$ws.Range("K132").Value = 16130.727
$ws.Range("L132").Value = 42755.781
$ws.Range("M132").Value = -13600.727
$ws.Range("N132").Value = -47815.781
$ws.Range("H138").Value = 4587.705
$ws.Range("I138").Value = 1133.0625
$ws.Range("K138").Value = 3399.1875
$ws.Range("M138").Value = 1740.8125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3405.3333
$ws.Range("I45").Value = 3086.6
$ws.Range("K45").Value = 3086.6
$ws.Range("M45").Value = -2709.6
$ws.Range("H61").Value = 11132.823
$ws.Range("I61").Value = 18586.111
$ws.Range("J61").Value = 2747.875
$ws.Range("K61").Value = 18586.111
$ws.Range("L61").Value = 2747.875
$ws.Range("M61").Value = -18374.111
$ws.Range("N61").Value = -3171.875
$ws.Range("H136").Value = 11132.823
$ws.Range("I136").Value = 18586.111
$ws.Range("J136").Value = 2747.875
$ws.Range("K136").Value = 55758.333
$ws.Range("L136").Value = 8243.625
$ws.Range("M136").Value = -53208.333
$ws.Range("N136").Value = -13343.625

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2037.1666
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2037.1666
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 2084773.4
$ws.Range("I99").Value = 2977319.2
$ws.Range("K99").Value = 2977319.2
$ws.Range("M99").Value = -2975821.2
$ws.Range("H107").Value = 1750.3334
$ws.Range("I107").Value = 1742.5294
$ws.Range("K107").Value = 1742.5294
$ws.Range("M107").Value = 177.4706000000001
$ws.Range("H134").Value = 4461.1113
$ws.Range("I134").Value = 1691.6666
$ws.Range("K134").Value = 5074.9998
$ws.Range("M134").Value = -2539.9998
$ws.Range("H138").Value = 91000
$ws.Range("J138").Value = 91000
$ws.Range("L138").Value = 91000
$ws.Range("N138").Value = -101280

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5832.891
$ws.Range("I31").Value = 2010.3158
$ws.Range("J31").Value = 8522.852000000001
$ws.Range("K31").Value = 2010.3158
$ws.Range("L31").Value = 8522.852000000001
$ws.Range("M31").Value = -1715.3158
$ws.Range("N31").Value = -9112.852000000001
$ws.Range("H34").Value = 5832.891
$ws.Range("I34").Value = 2010.3158
$ws.Range("J34").Value = 8522.852000000001
$ws.Range("K34").Value = 2010.3158
$ws.Range("L34").Value = 8522.852000000001
$ws.Range("M34").Value = -1808.3158
$ws.Range("N34").Value = -8926.852000000001
$ws.Range("H99").Value = 6609.5
$ws.Range("J99").Value = 9475
$ws.Range("L99").Value = 9475
$ws.Range("N99").Value = -12471
$ws.Range("H107").Value = 512.8570999999999
$ws.Range("I107").Value = 425.72726
$ws.Range("J107").Value = 832.3333
$ws.Range("K107").Value = 425.72726
$ws.Range("L107").Value = 832.3333
$ws.Range("M107").Value = 1494.27274
$ws.Range("N107").Value = -4672.3333
$ws.Range("H126").Value = 6609.5
$ws.Range("J126").Value = 9475
$ws.Range("L126").Value = 28425
$ws.Range("N126").Value = -33365

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 917.3333
$ws.Range("I5").Value = 757.75
$ws.Range("J5").Value = 949.25
$ws.Range("K5").Value = 2273.25
$ws.Range("L5").Value = 2847.75
$ws.Range("M5").Value = -2161.25
$ws.Range("N5").Value = -3071.75
$ws.Range("H64").Value = 2515983.8
$ws.Range("J64").Value = 2873981.5
$ws.Range("L64").Value = 8621944.5
$ws.Range("N64").Value = -8622484.5
$ws.Range("H67").Value = 2515983.8
$ws.Range("J67").Value = 2873981.5
$ws.Range("L67").Value = 8621944.5
$ws.Range("N67").Value = -8623816.5
$ws.Range("H107").Value = 1342.8667
$ws.Range("J107").Value = 1397.3572
$ws.Range("L107").Value = 4192.071599999999
$ws.Range("N107").Value = -8032.071599999999
$ws.Range("H113").Value = 762.5
$ws.Range("I113").Value = 666
$ws.Range("K113").Value = 1998
$ws.Range("M113").Value = 172
$ws.Range("H135").Value = 917.3333
$ws.Range("I135").Value = 757.75
$ws.Range("J135").Value = 949.25
$ws.Range("K135").Value = 6819.75
$ws.Range("L135").Value = 8543.25
$ws.Range("M135").Value = -4284.75
$ws.Range("N135").Value = -13613.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4611.615
$ws.Range("I122").Value = 2561.5
$ws.Range("K122").Value = 7684.5
$ws.Range("M122").Value = -5234.5
$ws.Range("H138").Value = 52214.5
$ws.Range("J138").Value = 52214.5
$ws.Range("L138").Value = 52214.5
$ws.Range("N138").Value = -62494.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9093028
$ws.Range("I16").Value = 11112886
$ws.Range("K16").Value = 11112886
$ws.Range("M16").Value = -11112716
$ws.Range("H40").Value = 25642210
$ws.Range("I40").Value = 1282.8334
$ws.Range("K40").Value = 1282.8334
$ws.Range("M40").Value = -1146.8334
$ws.Range("H46").Value = 7411.522
$ws.Range("J46").Value = 7411.522
$ws.Range("L46").Value = 7411.522
$ws.Range("N46").Value = -7787.522
$ws.Range("H134").Value = 82619.336
$ws.Range("J134").Value = 82619.336
$ws.Range("L134").Value = 82619.336
$ws.Range("N134").Value = -92759.336
$ws.Range("H135").Value = 133328
$ws.Range("J135").Value = 133328
$ws.Range("L135").Value = 133328
$ws.Range("N135").Value = -143468

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 36481.332
$ws.Range("I2").Value = 36481.332
$ws.Range("K2").Value = 36481.332
$ws.Range("M2").Value = -36369.332
$ws.Range("H39").Value = 49999
$ws.Range("I39").Value = 49999
$ws.Range("K39").Value = 49999
$ws.Range("M39").Value = -49586
$ws.Range("H45").Value = 42062.5
$ws.Range("J45").Value = 42062.5
$ws.Range("L45").Value = 42062.5
$ws.Range("N45").Value = -43044.5
$ws.Range("H62").Value = 26218.5
$ws.Range("I62").Value = 17437.5
$ws.Range("K62").Value = 17437.5
$ws.Range("M62").Value = -16813.5
$ws.Range("H65").Value = 26218.5
$ws.Range("I65").Value = 17437.5
$ws.Range("K65").Value = 87187.5
$ws.Range("M65").Value = -84067.5
$ws.Range("H107").Value = 2403.0417
$ws.Range("I107").Value = 2633.8
$ws.Range("K107").Value = 7901.400000000001
$ws.Range("M107").Value = -5981.400000000001
$ws.Range("H113").Value = 998.8333
$ws.Range("I113").Value = 1025.5714
$ws.Range("J113").Value = 961.4
$ws.Range("K113").Value = 3076.7142
$ws.Range("L113").Value = 2884.2
$ws.Range("M113").Value = -906.7142000000003
$ws.Range("N113").Value = -7224.2
$ws.Range("H132").Value = 39694324
$ws.Range("I132").Value = 6174284
$ws.Range("J132").Value = 100030400
$ws.Range("K132").Value = 18522852
$ws.Range("L132").Value = 300091200
$ws.Range("M132").Value = -18520322
$ws.Range("N132").Value = -300096260
$ws.Range("H136").Value = 8769.109
$ws.Range("I136").Value = 3672.037
$ws.Range("J136").Value = 11760.869
$ws.Range("K136").Value = 11016.111
$ws.Range("L136").Value = 35282.607
$ws.Range("M136").Value = -8466.110999999999
$ws.Range("N136").Value = -40382.607
$ws.Range("H140").Value = 90669.25
$ws.Range("J140").Value = 94095.664
$ws.Range("L140").Value = 94095.664
$ws.Range("N140").Value = -104455.664

